$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column L (empty helper column between the "M-label" cell and the
# "O:V" overlap-count block). This shifts the old M column to L and the old
# O:V block to N:U, matching the recorded edit.
$ws.Columns("L").Delete()

# Text updates (order matters: new shared strings get appended in the
# order they are first referenced, so write these in the same order as
# the target sharedStrings table was built).
$ws.Range("A17").Value = "nel datset MATCHED crei un identificativo di ogni gruppo e ne conti la numerosità"
$ws.Range("A18").Value = "nel dataset ORIGIN etichetti tutti i record id che sono nel dataset MATCHED, con la numerosità delle sue ripetizioni e l'identificativo del gruppo"
$ws.Range("A7").Value = "ORIGIN: dataset dei "

# Update the selection shown when the file is reopened.
$ws.Range("A1:C5").Select() | Out-Null
